$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General Iteration")

# --- verboseOutput (row 8): FALSE -> TRUE ---
$ws.Range("C8").Value = $true

# --- rngSeed (row 12): add value 4 in C12 ---
$ws.Range("C12").Value = 4

# --- cohesion (row 27): C27 0.1 -> 3000 ---
$ws.Range("C27").Value = 3000

# --- heightFactorPower (row 28): C28 0.5 -> 1 ---
$ws.Range("C28").Value = 1

# --- cohesionAscensionIgnore (row 29): C29 -3 -> 0.5, drop D29:J29 ---
$ws.Range("C29").Value = 0.5
$ws.Range("D29:J29").ClearContents()

# --- cohesionAscensionMax (row 30): C30 stays 10, drop D30:J30 ---
$ws.Range("D30:J30").ClearContents()

# --- cohPower (row 31): C31 0.5 -> 2 ---
$ws.Range("C31").Value = 2

# --- sepPower (row 32): C32 0.5 -> 0, drop D32 ---
$ws.Range("C32").Value = 0
$ws.Range("D32").ClearContents()

# --- alignment (row 34): C34 0.01 -> 0.1 ---
$ws.Range("C34").Value = 0.1

# --- alignmentHeightWidth / sepPower (row 36): C36 -6 -> -2 ---
$ws.Range("C36").Value = -2

# --- waggle (row 38): C38 1E-4 -> 0 ---
$ws.Range("C38").Value = 0

# --- New rows 98-100: Scores section ---
$ws.Range("A98").Value = "Scores"
$ws.Range("A98").Font.Bold = $true

$ws.Range("A99").Value = "How many boxes the map should be divided into"
$ws.Range("B99").Value = "mapDivResolution"
$ws.Range("C99").Value = 10

$ws.Range("A100").Value = "How often to check which divisions have been explored"
$ws.Range("B100").Value = "mapDivFrameSkip"
$ws.Range("C100").Value = 10

# --- Update selection / view to match the saved workbook state ---
[void]$ws.Range("B96").Select()
try {
    $excel.ActiveWindow.ScrollRow = 91
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # scroll position is a cosmetic viewport setting; ignore if unsupported
}
